$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force Text format on the data range so numeric-looking
# strings (e.g. "207.65", "12.40") are preserved exactly as text,
# matching the original inline-string cell contents.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

# Row 2
$ws.Range('D2').Value = '68.291.79'
$ws.Range('E2').Value = '  +1.54%  '

# Row 3
$ws.Range('D3').Value = '3.598.67'
$ws.Range('E3').Value = '  +0.30%  '

# Row 4
$ws.Range('E4').Value = '  -0.10%  '

# Row 5
$ws.Range('D5').Value = '207.65'
$ws.Range('E5').Value = '  +7.86%  '

# Row 6
$ws.Range('D6').Value = '570.11'
$ws.Range('E6').Value = '  -1.30%  '

# Row 7
$ws.Range('E7').Value = '  -0.62%  '

# Row 8
$ws.Range('E8').Value = '  -0.03%  '

# Row 9
$ws.Range('D9').Value = '0.685'
$ws.Range('E9').Value = '  +0.75%  '

# Row 10
$ws.Range('D10').Value = '63.99'
$ws.Range('E10').Value = '  +13.98%  '

# Row 11
$ws.Range('E11').Value = '  -1.30%  '

# Row 12
$ws.Range('E12').Value = '  +2.94%  '

# Row 13
$ws.Range('E13').Value = '  +5.11%  '

# Row 14
$ws.Range('D14').Value = '4.170.98'
$ws.Range('E14').Value = '  -0.10%  '

# Row 15
$ws.Range('D15').Value = '3.601.37'
$ws.Range('E15').Value = '  +0.39%  '

# Row 16
$ws.Range('D16').Value = '19.22'
$ws.Range('E16').Value = '  +4.43%  '

# Row 17
$ws.Range('E17').Value = '  +0.46%  '

# Row 18
$ws.Range('D18').Value = '68.087.31'
$ws.Range('E18').Value = '  +1.23%  '

# Row 19
$ws.Range('E19').Value = '  +0.43%  '

# Row 20
$ws.Range('E20').Value = '  +0.18%  '

# Row 21
$ws.Range('D21').Value = '405.38'
$ws.Range('E21').Value = '  +0.77%  '

# Row 22
$ws.Range('E22').Value = '  -1.02%  '

# Row 23
$ws.Range('D23').Value = '12.40'
$ws.Range('E23').Value = '  +8.59%  '

# Row 24
$ws.Range('D24').Value = '84.89'
$ws.Range('E24').Value = '  -1.18%  '

# Row 25
$ws.Range('D25').Value = '2.90'
$ws.Range('E25').Value = '  -1.33%  '

# Row 26
$ws.Range('D26').Value = '12.55'
$ws.Range('E26').Value = '  +0.48%  '

# Row 27
$ws.Range('D27').Value = '3.86'
$ws.Range('E27').Value = '  +5.79%  '

# Row 28
$ws.Range('D28').Value = '9.27'
$ws.Range('E28').Value = '  +3.31%  '

# Row 29
$ws.Range('D29').Value = '7.63'
$ws.Range('E29').Value = '  -0.40%  '

# Row 30
$ws.Range('D30').Value = '31.62'
$ws.Range('E30').Value = '  +1.21%  '

# Row 31
$ws.Range('D31').Value = '696.32'
$ws.Range('E31').Value = '  +9.82%  '

# Row 32
$ws.Range('D32').Value = '12.19'
$ws.Range('E32').Value = '  -0.01%  '

# Row 33
$ws.Range('E33').Value = '  -1.01%  '

# Row 34
$ws.Range('D34').Value = '63.66'
$ws.Range('E34').Value = '  -0.57%  '

# Row 35
$ws.Range('D35').Value = '41.66'
$ws.Range('E35').Value = '  -2.16%  '

# Row 36
$ws.Range('D36').Value = '0.412'
$ws.Range('E36').Value = '  +2.98%  '

# Row 37
$ws.Range('E37').Value = '  +0.27%  '

# Row 38
$ws.Range('E38').Value = '  +8.29%  '

# Row 39
$ws.Range('E39').Value = '  -2.95%  '

# Row 40
$ws.Range('D40').Value = '3.19'
$ws.Range('E40').Value = '  +21.14%  '

# Row 41
$ws.Range('D41').Value = '3.172.61'
$ws.Range('E41').Value = '  -0.09%  '

# Row 42
$ws.Range('E42').Value = '  -0.68%  '

# Row 43
$ws.Range('D43').Value = '0.998'
$ws.Range('E43').Value = '  -0.32%  '

# Row 44
$ws.Range('E44').Value = '  -1.68%  '

# Row 45
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').Value = '0.0414'
$ws.Range('E45').Value = '  -0.60%  '

# Row 46
$ws.Range('B46').Value = 'WEMIXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D46').Value = '2.76'
$ws.Range('E46').Value = '  +8.61%  '

# Row 47
$ws.Range('D47').Value = '3.12'
$ws.Range('E47').Value = '  -0.09%  '

# Row 48
$ws.Range('E48').Value = '  +0.42%  '

# Row 49
$ws.Range('D49').Value = '8.79'
$ws.Range('E49').Value = '  +2.23%  '

# Row 50
$ws.Range('D50').Value = '139.15'
$ws.Range('E50').Value = '  -1.84%  '

# Row 51
$ws.Range('E51').Value = '  -1.18%  '

# Restore the default (Normal) style on the data range so the
# underlying cell formatting matches the original workbook.
$dataRange.Style = "Normal"